# r49-73.py: Create script -- finish hospital categorization later
#
# Fills in September 2020 ("sep_2020", column G) figures for the hospital /
# disposition-categorization rows (County Hospital, State Hospital, the
# Psychiatric Hospitalizations subtotal, and the various "Dispositioned To"
# community categories), and carries each value into the "SFY 2021 Total"
# column (Q) since September was the only month with data so far.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> value entered for sep_2020 (column G); SFY 2021 Total (column Q)
# picks up the same figure since it is the only non-zero month reported.
$updates = [ordered]@{
    52 = 2   # County Hospital
    53 = 1   # State Hospital
    54 = 3   # Total - Psychiatric Hospitalizations
    55 = 1   # Legal System Involved (Jail, Police, Court)
    56 = 4   # Affiliated Emergency Service
    59 = 2   # Partial Care/Partial Hospitalization
    61 = 5   # ICMS
    66 = 1   # Other Mental Health Services (e.g. private practicioner)
    71 = 1   # INVOLUNTARY OUT-PATIENT COMMITMENT
    72 = 21  # Other
    73 = 25  # None
    74 = 60  # Total unduplicated DSC consumers discharged to community
}

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    $ws.Cells.Item($row, 7).Value = $value   # column G = sep_2020
    $ws.Cells.Item($row, 17).Value = $value  # column Q = SFY 2021 Total
}

# Restore the reviewer's scroll/selection state from the saved view.
$ws.Activate()
$ws.Range("A67").Select()
$ws.Range("I71").Select()
